$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for rows 2..70 (UF, Ano, Taxa de atendimento).
# A row for 01/01/2023 was added for each UF (Brasil, Nordeste, Sergipe),
# shifting the existing per-UF blocks down by one row, and several
# trailing-year figures were refreshed with updated source values.
$data = @(
    @("Brasil","01/01/2001",0.3556014597832986),
    @("Brasil","01/01/2002",0.4091636330542122),
    @("Brasil","01/01/2003",0.5044635381561404),
    @("Brasil","01/01/2004",0.5309432578199587),
    @("Brasil","01/01/2005",0.4890586980760805),
    @("Brasil","01/01/2006",0.4620897956711576),
    @("Brasil","01/01/2007",0.4966081127698825),
    @("Brasil","01/01/2008",0.4599472186023644),
    @("Brasil","01/01/2009",0.5105518385462808),
    @("Brasil","01/01/2010",0.6136231887283073),
    @("Brasil","01/01/2011",0.4826639242056643),
    @("Brasil","01/01/2012",1.528346348537724),
    @("Brasil","01/01/2013",1.661613778404539),
    @("Brasil","01/01/2014",1.597383886264928),
    @("Brasil","01/01/2015",1.240272622756693),
    @("Brasil","01/01/2016",1.29005466803939),
    @("Brasil","01/01/2017",1.24313952965738),
    @("Brasil","01/01/2018",1.204507731023762),
    @("Brasil","01/01/2019",1.136973218646514),
    @("Brasil","01/01/2020",0.9804530571292559),
    @("Brasil","01/01/2021",1.065176199645811),
    @("Brasil","01/01/2022",1.087798046816562),
    @("Brasil","01/01/2023",0.3004981230324768),
    @("Nordeste","01/01/2001",0.4772744654469883),
    @("Nordeste","01/01/2002",0.5532575945051468),
    @("Nordeste","01/01/2003",0.6204436706082329),
    @("Nordeste","01/01/2004",0.6269770215579353),
    @("Nordeste","01/01/2005",0.5949981691616898),
    @("Nordeste","01/01/2006",0.5419991935493135),
    @("Nordeste","01/01/2007",0.5858672095479524),
    @("Nordeste","01/01/2008",0.4859552888009054),
    @("Nordeste","01/01/2009",0.5531732981338393),
    @("Nordeste","01/01/2010",0.6577759822211939),
    @("Nordeste","01/01/2011",0.5229696404338543),
    @("Nordeste","01/01/2012",1.874680179165634),
    @("Nordeste","01/01/2013",2.133323612104753),
    @("Nordeste","01/01/2014",2.124534468979905),
    @("Nordeste","01/01/2015",1.444567115074851),
    @("Nordeste","01/01/2016",1.511515509064918),
    @("Nordeste","01/01/2017",1.44488291740586),
    @("Nordeste","01/01/2018",1.374213316911246),
    @("Nordeste","01/01/2019",1.256496275158019),
    @("Nordeste","01/01/2020",1.06910676406558),
    @("Nordeste","01/01/2021",1.166474645942877),
    @("Nordeste","01/01/2022",1.207480173827043),
    @("Nordeste","01/01/2023",0.322157260470111),
    @("Sergipe","01/01/2001",0.2340698640198239),
    @("Sergipe","01/01/2002",0.39975420518465),
    @("Sergipe","01/01/2003",0.5628238604846239),
    @("Sergipe","01/01/2004",0.4621264110910339),
    @("Sergipe","01/01/2005",0.5388747331451003),
    @("Sergipe","01/01/2006",0.4420086778692693),
    @("Sergipe","01/01/2007",0.4434061928072969),
    @("Sergipe","01/01/2008",0.3467491695790255),
    @("Sergipe","01/01/2009",0.4085367151094751),
    @("Sergipe","01/01/2010",0.4532548490281976),
    @("Sergipe","01/01/2011",0.3123523449870334),
    @("Sergipe","01/01/2012",1.385008639330851),
    @("Sergipe","01/01/2013",1.492646399277098),
    @("Sergipe","01/01/2014",1.346030902059082),
    @("Sergipe","01/01/2015",0.9489994144273937),
    @("Sergipe","01/01/2016",0.8865412141247279),
    @("Sergipe","01/01/2017",0.7644810349897089),
    @("Sergipe","01/01/2018",0.7216124650668998),
    @("Sergipe","01/01/2019",0.6650468019211764),
    @("Sergipe","01/01/2020",0.7519397426027848),
    @("Sergipe","01/01/2021",0.8725628497491614),
    @("Sergipe","01/01/2022",0.9812372147795632),
    @("Sergipe","01/01/2023",0.3053398335988244)
)

$rowCount = $data.Count
$lastRow = 1 + $rowCount

$bRange = $ws.Range($ws.Cells.Item(2,2), $ws.Cells.Item($lastRow,2))

# Column B holds dd/mm/yyyy-look-alike text ("01/01/2001", ...). Force
# text format first so Excel does not silently coerce the strings into
# date serial numbers, then write, then drop back to the default style
# so no stray number format lingers on the cells.
$bRange.NumberFormat = "@"

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

$bRange.Style = "Normal"

Write-Host "Wrote $rowCount data rows (A2:C$lastRow)"
